# Add the "Lore of Beasts" spells to the Spells sheet (rows 84-91),
# then leave the Spells tab as the active/selected sheet with the same
# cell selection the author ended up with.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spells")

# Row 84 - Amber Spear
$ws.Range("A84").Value = "Amber Spear"
$ws.Range("B84").Value = "New"
$ws.Range("C84").Value = 3
$ws.Range("D84").Value = "Evocation"
$ws.Range("E84").Value = "No"
$ws.Range("F84").Value = "No"
$ws.Range("G84").Value = "No"
$ws.Range("H84").Value = "Yes"
$ws.Range("I84").Value = "No"
$ws.Range("J84").Value = "Yes"
$ws.Range("K84").Value = "No"
$ws.Range("L84").Value = "No"
$ws.Range("M84").Value = "No"
$ws.Range("N84").Value = "1.0.0"
$ws.Range("O84").Value = "Complete"
$ws.Range("P84").Value = "Publicly Released"
$ws.Range("Q84").Value = "Not on website"

# Row 85 - Aspect of the Beast
$ws.Range("A85").Value = "Aspect of the Beast"
$ws.Range("B85").Value = "New"
$ws.Range("C85").Value = 4
$ws.Range("D85").Value = "Trasmutation"
$ws.Range("E85").Value = "No"
$ws.Range("F85").Value = "No"
$ws.Range("G85").Value = "No"
$ws.Range("H85").Value = "Yes"
$ws.Range("I85").Value = "No"
$ws.Range("J85").Value = "Yes"
$ws.Range("K85").Value = "No"
$ws.Range("L85").Value = "No"
$ws.Range("M85").Value = "No"
$ws.Range("N85").Value = "1.0.0"
$ws.Range("O85").Value = "Complete"
$ws.Range("P85").Value = "Publicly Released"
$ws.Range("Q85").Value = "Not on website"

# Row 86 - Bestial Spirit
$ws.Range("A86").Value = "Bestial Spirit"
$ws.Range("B86").Value = "New"
$ws.Range("C86").Value = 4
$ws.Range("D86").Value = "Conjuration"
$ws.Range("E86").Value = "No"
$ws.Range("F86").Value = "No"
$ws.Range("G86").Value = "No"
$ws.Range("H86").Value = "Yes"
$ws.Range("I86").Value = "No"
$ws.Range("J86").Value = "Yes"
$ws.Range("K86").Value = "No"
$ws.Range("L86").Value = "No"
$ws.Range("M86").Value = "No"
$ws.Range("N86").Value = "1.0.0"
$ws.Range("O86").Value = "Complete"
$ws.Range("P86").Value = "Publicly Released"
$ws.Range("Q86").Value = "Not on website"

# Row 87 - Call of the Pack
$ws.Range("A87").Value = "Call of the Pack"
$ws.Range("B87").Value = "New"
$ws.Range("C87").Value = 5
$ws.Range("D87").Value = "Enchantment"
$ws.Range("E87").Value = "No"
$ws.Range("F87").Value = "No"
$ws.Range("G87").Value = "No"
$ws.Range("H87").Value = "No"
$ws.Range("I87").Value = "No"
$ws.Range("J87").Value = "Yes"
$ws.Range("K87").Value = "No"
$ws.Range("L87").Value = "No"
$ws.Range("M87").Value = "No"
$ws.Range("N87").Value = "1.0.0"
$ws.Range("O87").Value = "Complete"
$ws.Range("P87").Value = "Publicly Released"
$ws.Range("Q87").Value = "Not on website"

# Row 88 - Impenetrable Pelt
$ws.Range("A88").Value = "Impenetrable Pelt"
$ws.Range("B88").Value = "New"
$ws.Range("C88").Value = 5
$ws.Range("D88").Value = "Trasmutation"
$ws.Range("E88").Value = "No"
$ws.Range("F88").Value = "No"
$ws.Range("G88").Value = "No"
$ws.Range("H88").Value = "Yes"
$ws.Range("I88").Value = "No"
$ws.Range("J88").Value = "Yes"
$ws.Range("K88").Value = "No"
$ws.Range("L88").Value = "No"
$ws.Range("M88").Value = "No"
$ws.Range("N88").Value = "1.0.0"
$ws.Range("O88").Value = "Complete"
$ws.Range("P88").Value = "Publicly Released"
$ws.Range("Q88").Value = "Not on website"

# Row 89 - Monstrous Transformation
$ws.Range("A89").Value = "Monstrous Transformation"
$ws.Range("B89").Value = "New"
$ws.Range("C89").Value = 6
$ws.Range("D89").Value = "Trasmutation"
$ws.Range("E89").Value = "No"
$ws.Range("F89").Value = "No"
$ws.Range("G89").Value = "No"
$ws.Range("H89").Value = "Yes"
$ws.Range("I89").Value = "No"
$ws.Range("J89").Value = "No"
$ws.Range("K89").Value = "Yes"
$ws.Range("L89").Value = "Yes"
$ws.Range("M89").Value = "No"
$ws.Range("N89").Value = "1.0.0"
$ws.Range("O89").Value = "Complete"
$ws.Range("P89").Value = "Publicly Released"
$ws.Range("Q89").Value = "Not on website"

# Row 90 - Primal Dominance
$ws.Range("A90").Value = "Primal Dominance"
$ws.Range("B90").Value = "New"
$ws.Range("C90").Value = 2
$ws.Range("D90").Value = "Enchantment"
$ws.Range("E90").Value = "No"
$ws.Range("F90").Value = "Yes"
$ws.Range("G90").Value = "No"
$ws.Range("H90").Value = "Yes"
$ws.Range("I90").Value = "No"
$ws.Range("J90").Value = "Yes"
$ws.Range("K90").Value = "No"
$ws.Range("L90").Value = "No"
$ws.Range("M90").Value = "No"
$ws.Range("N90").Value = "1.0.0"
$ws.Range("O90").Value = "Complete"
$ws.Range("P90").Value = "Publicly Released"
$ws.Range("Q90").Value = "Not on website"

# Row 91 - Summon Flock
$ws.Range("A91").Value = "Summon Flock"
$ws.Range("B91").Value = "New"
$ws.Range("C91").Value = 1
$ws.Range("D91").Value = "Conjuration"
$ws.Range("E91").Value = "No"
$ws.Range("F91").Value = "No"
$ws.Range("G91").Value = "No"
$ws.Range("H91").Value = "Yes"
$ws.Range("I91").Value = "No"
$ws.Range("J91").Value = "Yes"
$ws.Range("K91").Value = "No"
$ws.Range("L91").Value = "No"
$ws.Range("M91").Value = "No"
$ws.Range("N91").Value = "1.0.0"
$ws.Range("O91").Value = "Complete"
$ws.Range("P91").Value = "Publicly Released"
$ws.Range("Q91").Value = "Not on website"

# Make Spells the active/selected sheet (moves tabSelected off whichever
# sheet had it before) and leave the same cell selected as in the author's
# final view.
$ws.Select()
$ws.Range("E87").Select()
